# Updates the cryptos price table (columns B-E, rows 2-51) to the latest
# scraped snapshot, as produced by the scheduled "Updated cryptos list"
# GitHub Actions job.
#
# Notes:
#  - Column D ("Price") values are stored as literal text in the workbook
#    (e.g. "27.291.14", "1.000", "0.000008655") rather than numbers, since
#    they use a thousands-separator "." convention that would otherwise be
#    misinterpreted. We prefix those assignments with a leading apostrophe
#    (the same trick a person typing into Excel would use) so the COM
#    layer stores them as text instead of silently parsing them as numbers.
#  - Column E ("Volume(1h)") values already contain padding spaces, so
#    Excel naturally keeps them as text without any extra handling.
#  - A few coins swapped rank/position between rows (13-15, 24-25), which
#    shows up as updates to the Coin name (B) and Link (C) cells as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.291.14"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "'1.903.33"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'306.49"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").Value = "'0.5408"
$ws.Range("E7").Value = "  +3.81%  "

$ws.Range("E8").Value = "  +1.21%  "

$ws.Range("D9").Value = "'0.07299"
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").Value = "'22.07"
$ws.Range("E10").Value = "  +4.22%  "

$ws.Range("D11").Value = "'0.9011"
$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("D12").Value = "'0.08187"
$ws.Range("E12").Value = "  -0.26%  "

# Rows 13-15: Litecoin / Polkadot / WrappedEther rotated positions.
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.914.37"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'95.50"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.363"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("E17").Value = "  +1.91%  "

$ws.Range("D18").Value = "'0.000008655"
$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("D20").Value = "'27.320.49"
$ws.Range("E20").Value = "  +0.60%  "

$ws.Range("D21").Value = "'5.050"
$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("D22").Value = "'10.84"
$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("D23").Value = "'6.509"
$ws.Range("E23").Value = "  +1.38%  "

# Rows 24-25: LidoDAOToken / Monero swapped positions.
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'148.85"
$ws.Range("E24").Value = "  +0.22%  "

$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.307"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("D27").Value = "'1.758"
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("D28").Value = "'116.51"
$ws.Range("E28").Value = "  +1.00%  "

$ws.Range("D29").Value = "'4.837"
$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("D30").Value = "'4.662"
$ws.Range("E30").Value = "  -4.05%  "

$ws.Range("D31").Value = "'0.09169"
$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("D32").Value = "'0.8249"
$ws.Range("E32").Value = "  +3.81%  "

$ws.Range("D33").Value = "'0.05067"
$ws.Range("E33").Value = "  +0.83%  "

$ws.Range("D34").Value = "'1.224"
$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").Value = "'3.022"
$ws.Range("E35").Value = "  +1.82%  "

$ws.Range("E36").Value = "  -3.49%  "

$ws.Range("D37").Value = "'2.679"
$ws.Range("E37").Value = "  +2.47%  "

$ws.Range("D38").Value = "'0.5991"
$ws.Range("E38").Value = "  +4.66%  "

$ws.Range("D39").Value = "'0.01995"
$ws.Range("E39").Value = "  -0.33%  "

$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("D41").Value = "'9.274"
$ws.Range("E41").Value = "  +2.92%  "

$ws.Range("D42").Value = "'6.668"

$ws.Range("D43").Value = "'115.86"
$ws.Range("E43").Value = "  -0.47%  "

$ws.Range("D44").Value = "'0.5118"
$ws.Range("E44").Value = "  +4.97%  "

$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("D46").Value = "'10.23"
$ws.Range("E46").Value = "  +0.93%  "

$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").Value = "'1.636"
$ws.Range("E48").Value = "  +0.93%  "

$ws.Range("D49").Value = "'38.08"
$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("D50").Value = "'0.06090"
$ws.Range("E50").Value = "  +2.76%  "

$ws.Range("D51").Value = "'63.27"
$ws.Range("E51").Value = "  -0.68%  "
